# Auto-generated edit script: update Leve market-price columns (H-N) across all job sheets
# per scheduled Universalis price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 6499.75
$ws.Cells.Item(98, 9).Value = 5333
$ws.Cells.Item(98, 10).Value = 10000
$ws.Cells.Item(98, 11).Value = 5333
$ws.Cells.Item(98, 12).Value = 10000
$ws.Cells.Item(98, 13).Value = -3835
$ws.Cells.Item(98, 14).Value = -12996
$ws.Cells.Item(113, 8).Value = 8284.621999999999
$ws.Cells.Item(113, 9).Value = 5102.25
$ws.Cells.Item(113, 11).Value = 5102.25
$ws.Cells.Item(113, 13).Value = -1848.25
$ws.Cells.Item(122, 8).Value = 6499.75
$ws.Cells.Item(122, 9).Value = 5333
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 15999
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -13549
$ws.Cells.Item(122, 14).Value = -34900
$ws.Cells.Item(131, 8).Value = 12013.19
$ws.Cells.Item(131, 10).Value = 30728.715
$ws.Cells.Item(131, 12).Value = 92186.145
$ws.Cells.Item(131, 14).Value = -102266.145
$ws.Cells.Item(137, 8).Value = 1681.4651
$ws.Cells.Item(137, 10).Value = 1729.4546
$ws.Cells.Item(137, 12).Value = 5188.3638
$ws.Cells.Item(137, 14).Value = -10288.3638
$ws.Cells.Item(138, 8).Value = 8775063
$ws.Cells.Item(138, 9).Value = 1606.2307
$ws.Cells.Item(138, 10).Value = 11367221
$ws.Cells.Item(138, 11).Value = 4818.6921
$ws.Cells.Item(138, 12).Value = 34101663
$ws.Cells.Item(138, 13).Value = 321.3078999999998
$ws.Cells.Item(138, 14).Value = -34111943
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 486.17392
$ws.Cells.Item(2, 9).Value = 339.83334
$ws.Cells.Item(2, 11).Value = 339.83334
$ws.Cells.Item(2, 13).Value = -226.83334
$ws.Cells.Item(45, 8).Value = 5002.25
$ws.Cells.Item(45, 9).Value = 5878.0454
$ws.Cells.Item(45, 11).Value = 5878.0454
$ws.Cells.Item(45, 13).Value = -5501.0454
$ws.Cells.Item(61, 8).Value = 263746.2
$ws.Cells.Item(61, 9).Value = 3999.5
$ws.Cells.Item(61, 11).Value = 3999.5
$ws.Cells.Item(61, 13).Value = -3787.5
$ws.Cells.Item(104, 8).Value = 95914.3
$ws.Cells.Item(104, 10).Value = 95914.3
$ws.Cells.Item(104, 12).Value = 95914.3
$ws.Cells.Item(104, 14).Value = -102902.3
$ws.Cells.Item(116, 8).Value = 486.17392
$ws.Cells.Item(116, 9).Value = 339.83334
$ws.Cells.Item(116, 11).Value = 339.83334
$ws.Cells.Item(116, 13).Value = 1954.16666
$ws.Cells.Item(132, 8).Value = 4501.684
$ws.Cells.Item(132, 9).Value = 3908.7856
$ws.Cells.Item(132, 10).Value = 6161.8
$ws.Cells.Item(132, 11).Value = 11726.3568
$ws.Cells.Item(132, 12).Value = 18485.4
$ws.Cells.Item(132, 13).Value = -9196.356800000001
$ws.Cells.Item(132, 14).Value = -23545.4
$ws.Cells.Item(136, 8).Value = 263746.2
$ws.Cells.Item(136, 9).Value = 3999.5
$ws.Cells.Item(136, 11).Value = 11998.5
$ws.Cells.Item(136, 13).Value = -9448.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 486.17392
$ws.Cells.Item(3, 9).Value = 339.83334
$ws.Cells.Item(3, 11).Value = 339.83334
$ws.Cells.Item(3, 13).Value = -225.83334
$ws.Cells.Item(86, 8).Value = 3287.9736
$ws.Cells.Item(86, 9).Value = 3843.926
$ws.Cells.Item(86, 11).Value = 3843.926
$ws.Cells.Item(86, 13).Value = -2720.926
$ws.Cells.Item(89, 8).Value = 3287.9736
$ws.Cells.Item(89, 9).Value = 3843.926
$ws.Cells.Item(89, 11).Value = 19219.63
$ws.Cells.Item(89, 13).Value = -13603.63
$ws.Cells.Item(134, 8).Value = 2442.5293
$ws.Cells.Item(134, 9).Value = 2264.5625
$ws.Cells.Item(134, 10).Value = 5290
$ws.Cells.Item(134, 11).Value = 6793.6875
$ws.Cells.Item(134, 12).Value = 15870
$ws.Cells.Item(134, 13).Value = -4258.6875
$ws.Cells.Item(134, 14).Value = -20940
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 34700.555
$ws.Cells.Item(43, 10).Value = 34700.555
$ws.Cells.Item(43, 12).Value = 34700.555
$ws.Cells.Item(43, 14).Value = -35068.555
$ws.Cells.Item(53, 8).Value = 7299.6665
$ws.Cells.Item(53, 10).Value = 7299.6665
$ws.Cells.Item(53, 12).Value = 7299.6665
$ws.Cells.Item(53, 14).Value = -8513.666499999999
$ws.Cells.Item(58, 8).Value = 2695.95
$ws.Cells.Item(58, 9).Value = 1795.0714
$ws.Cells.Item(58, 11).Value = 1795.0714
$ws.Cells.Item(58, 13).Value = -1592.0714
$ws.Cells.Item(99, 8).Value = 2906.6365
$ws.Cells.Item(99, 9).Value = 2748
$ws.Cells.Item(99, 11).Value = 2748
$ws.Cells.Item(99, 13).Value = -1250
$ws.Cells.Item(101, 8).Value = 34700.555
$ws.Cells.Item(101, 10).Value = 34700.555
$ws.Cells.Item(101, 12).Value = 34700.555
$ws.Cells.Item(101, 14).Value = -41190.555
$ws.Cells.Item(124, 8).Value = 48647.5
$ws.Cells.Item(124, 10).Value = 48647.5
$ws.Cells.Item(124, 12).Value = 48647.5
$ws.Cells.Item(124, 14).Value = -53557.5
$ws.Cells.Item(126, 8).Value = 2906.6365
$ws.Cells.Item(126, 9).Value = 2748
$ws.Cells.Item(126, 11).Value = 8244
$ws.Cells.Item(126, 13).Value = -5774
$ws.Cells.Item(136, 8).Value = 2695.95
$ws.Cells.Item(136, 9).Value = 1795.0714
$ws.Cells.Item(136, 11).Value = 5385.2142
$ws.Cells.Item(136, 13).Value = -2835.2142
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value = 45003
$ws.Cells.Item(124, 9).Value = 30
$ws.Cells.Item(124, 11).Value = 90
$ws.Cells.Item(124, 13).Value = 4820
$ws.Cells.Item(126, 8).Value = 7300
$ws.Cells.Item(126, 9).Value = 7300
$ws.Cells.Item(126, 11).Value = 21900
$ws.Cells.Item(126, 13).Value = -16960
$ws.Cells.Item(138, 8).Value = 35723620
$ws.Cells.Item(138, 9).Value = 41675056
$ws.Cells.Item(138, 11).Value = 125025168
$ws.Cells.Item(138, 13).Value = -125020028
$ws.Cells.Item(140, 8).Value = 1696.4062
$ws.Cells.Item(140, 9).Value = 1696.4062
$ws.Cells.Item(140, 11).Value = 5089.2186
$ws.Cells.Item(140, 13).Value = 90.78139999999985
$ws.Cells.Item(141, 8).Value = 93839.27
$ws.Cells.Item(141, 9).Value = 1024.875
$ws.Cells.Item(141, 11).Value = 3074.625
$ws.Cells.Item(141, 13).Value = 2105.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 6622.6523
$ws.Cells.Item(132, 9).Value = 2683
$ws.Cells.Item(132, 11).Value = 8049
$ws.Cells.Item(132, 13).Value = -5519
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2193.122
$ws.Cells.Item(16, 9).Value = 1842.9412
$ws.Cells.Item(16, 10).Value = 3894
$ws.Cells.Item(16, 11).Value = 1842.9412
$ws.Cells.Item(16, 12).Value = 3894
$ws.Cells.Item(16, 13).Value = -1672.9412
$ws.Cells.Item(16, 14).Value = -4234
$ws.Cells.Item(22, 8).Value = 3142.7856
$ws.Cells.Item(22, 10).Value = 3258.3333
$ws.Cells.Item(22, 12).Value = 3258.3333
$ws.Cells.Item(22, 14).Value = -3848.3333
$ws.Cells.Item(27, 8).Value = 3142.7856
$ws.Cells.Item(27, 10).Value = 3258.3333
$ws.Cells.Item(27, 12).Value = 3258.3333
$ws.Cells.Item(27, 14).Value = -3472.3333
$ws.Cells.Item(40, 8).Value = 4423.357
$ws.Cells.Item(40, 9).Value = 3668.1428
$ws.Cells.Item(40, 10).Value = 5178.5713
$ws.Cells.Item(40, 11).Value = 3668.1428
$ws.Cells.Item(40, 12).Value = 5178.5713
$ws.Cells.Item(40, 13).Value = -3532.1428
$ws.Cells.Item(40, 14).Value = -5450.5713
$ws.Cells.Item(68, 8).Value = 2874.8333
$ws.Cells.Item(68, 10).Value = 2750
$ws.Cells.Item(68, 12).Value = 2750
$ws.Cells.Item(68, 14).Value = -4248
$ws.Cells.Item(71, 8).Value = 2874.8333
$ws.Cells.Item(71, 10).Value = 2750
$ws.Cells.Item(71, 12).Value = 13750
$ws.Cells.Item(71, 14).Value = -21238
$ws.Cells.Item(93, 8).Value = 1544.3939
$ws.Cells.Item(93, 9).Value = 1540.5358
$ws.Cells.Item(93, 11).Value = 1540.5358
$ws.Cells.Item(93, 13).Value = -292.5358000000001
$ws.Cells.Item(100, 8).Value = 401600.6
$ws.Cells.Item(100, 9).Value = 1751.5
$ws.Cells.Item(100, 10).Value = 668166.7
$ws.Cells.Item(100, 11).Value = 1751.5
$ws.Cells.Item(100, 12).Value = 668166.7
$ws.Cells.Item(100, 13).Value = -1210.5
$ws.Cells.Item(100, 14).Value = -669248.7
$ws.Cells.Item(101, 8).Value = 21078.6
$ws.Cells.Item(101, 10).Value = 21078.6
$ws.Cells.Item(101, 12).Value = 21078.6
$ws.Cells.Item(101, 14).Value = -27568.6
$ws.Cells.Item(134, 8).Value = 41237.375
$ws.Cells.Item(134, 10).Value = 41237.375
$ws.Cells.Item(134, 12).Value = 41237.375
$ws.Cells.Item(134, 14).Value = -51377.375
$ws.Cells.Item(136, 8).Value = 3158.5757
$ws.Cells.Item(136, 9).Value = 2782.7778
$ws.Cells.Item(136, 10).Value = 4849.6665
$ws.Cells.Item(136, 11).Value = 8348.3334
$ws.Cells.Item(136, 12).Value = 14548.9995
$ws.Cells.Item(136, 13).Value = -5798.3334
$ws.Cells.Item(136, 14).Value = -19648.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4670.625
$ws.Cells.Item(62, 9).Value = 2675
$ws.Cells.Item(62, 10).Value = 7996.6665
$ws.Cells.Item(62, 11).Value = 2675
$ws.Cells.Item(62, 12).Value = 7996.6665
$ws.Cells.Item(62, 13).Value = -2051
$ws.Cells.Item(62, 14).Value = -9244.666499999999
$ws.Cells.Item(65, 8).Value = 4670.625
$ws.Cells.Item(65, 9).Value = 2675
$ws.Cells.Item(65, 10).Value = 7996.6665
$ws.Cells.Item(65, 11).Value = 13375
$ws.Cells.Item(65, 12).Value = 39983.3325
$ws.Cells.Item(65, 13).Value = -10255
$ws.Cells.Item(65, 14).Value = -46223.3325
$ws.Cells.Item(81, 8).Value = 10140.375
$ws.Cells.Item(81, 9).Value = 22049.8
$ws.Cells.Item(81, 10).Value = 4727
$ws.Cells.Item(81, 11).Value = 44099.6
$ws.Cells.Item(81, 12).Value = 9454
$ws.Cells.Item(81, 13).Value = -43038.6
$ws.Cells.Item(81, 14).Value = -11576
$ws.Cells.Item(84, 8).Value = 10140.375
$ws.Cells.Item(84, 9).Value = 22049.8
$ws.Cells.Item(84, 10).Value = 4727
$ws.Cells.Item(84, 11).Value = 220498
$ws.Cells.Item(84, 12).Value = 47270
$ws.Cells.Item(84, 13).Value = -215194
$ws.Cells.Item(84, 14).Value = -57878
$ws.Cells.Item(105, 8).Value = 55410.4
$ws.Cells.Item(105, 10).Value = 55410.4
$ws.Cells.Item(105, 12).Value = 55410.4
$ws.Cells.Item(105, 14).Value = -62398.4
$ws.Cells.Item(136, 8).Value = 2239.2424
$ws.Cells.Item(136, 9).Value = 1358.375
$ws.Cells.Item(136, 11).Value = 4075.125
$ws.Cells.Item(136, 13).Value = -1525.125
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()
